$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.1121418643933389
$ws.Range("H2").Value = -14.66958496187364
$ws.Range("I2").Value = 7.354661273015827
$ws.Range("G3").Value = 0.1698111566353668
$ws.Range("H3").Value = 90.71934393713722
$ws.Range("G4").Value = -0.6541213138540207
$ws.Range("H4").Value = -6.301826861582926
$ws.Range("G5").Value = -0.648846692997273
$ws.Range("H5").Value = -6.264964206256273
$ws.Range("G6").Value = 0.149913777806357
$ws.Range("H6").Value = -39.08612898743718
$ws.Range("G7").Value = 0.3285302406496415
$ws.Range("H7").Value = 100.5410478824704
$ws.Range("G8").Value = 0.1186908439128481
$ws.Range("H8").Value = -28.20058476418997
$ws.Range("G9").Value = 0.2177951649834682
$ws.Range("H9").Value = 11.63684998636027
$ws.Range("G10").Value = -0.1289618669668079
$ws.Range("H10").Value = -125.66948249397
$ws.Range("G11").Value = -0.1110998533787646
$ws.Range("H11").Value = 6.456662387099604
$ws.Range("G12").Value = 0.1912955552428952
$ws.Range("H12").Value = 20.28951415318874
$ws.Range("G13").Value = 0.256380997536664
$ws.Range("H13").Value = 24.66148951006526
$ws.Range("G14").Value = 0.230356093912151
$ws.Range("H14").Value = 21.64141724881054
$ws.Range("G15").Value = 0.2276018980975815
$ws.Range("H15").Value = -8.924784317773394
$ws.Range("G16").Value = 0.006153625119849877
$ws.Range("H16").Value = -83.13162425233878
$ws.Range("G17").Value = 0.005021179681234354
$ws.Range("H17").Value = -85.84410444127727
$ws.Range("G18").Value = 0.08340672811650122
$ws.Range("H18").Value = -51.87542805717521
$ws.Range("G19").Value = 0.167108847032178
$ws.Range("H19").Value = 32.8866282652608
$ws.Range("G20").Value = 0.07762632544126712
$ws.Range("H20").Value = -32.29396356331548
$ws.Range("G21").Value = 0.1377111428131169
$ws.Range("H21").Value = 37.16409681289914
$ws.Range("G22").Value = 0.10591232615147
$ws.Range("H22").Value = 12.43689389173884
$ws.Range("G23").Value = 0.07934436199580493
$ws.Range("H23").Value = -26.86480701531849
$ws.Range("G24").Value = -0.1525691624515212
$ws.Range("H24").Value = -22.42323883724434
$ws.Range("G25").Value = -0.2148527433246224
$ws.Range("H25").Value = 3.417039223982303
$ws.Range("G26").Value = 0.1738512090749327
$ws.Range("H26").Value = 9.352888008466881
$ws.Range("G27").Value = 0.2074483805900574
$ws.Range("H27").Value = 3.482901070135082
$ws.Range("G28").Value = 0.03557104591645446
$ws.Range("H28").Value = 542.9520759543811
$ws.Range("G29").Value = -0.01232626068449046
$ws.Range("H29").Value = -180.1578099278001
